# Update the cached "datetimeFigureOut" field text (10/24/22 -> 10/26/22)
# on every Date placeholder: the slide master, each of the 11 slide
# layouts, and the notes master.

$p = $ppt.ActivePresentation

# --- Slide Master: Date Placeholder 3 (shape index 3) ---
$master = $p.SlideMaster
$master.Shapes.Item(3).TextFrame.TextRange.Text = "10/26/22"

# --- Slide Layouts: Date Placeholder shape index varies per layout ---
$layoutDateShapeIndex = @{
    1 = 3
    2 = 3
    3 = 3
    4 = 4
    5 = 6
    6 = 2
    7 = 1
    8 = 4
    9 = 4
    10 = 3
    11 = 3
}

foreach ($i in 1..11) {
    $layout = $master.CustomLayouts.Item($i)
    $shapeIdx = $layoutDateShapeIndex[$i]
    $layout.Shapes.Item($shapeIdx).TextFrame.TextRange.Text = "10/26/22"
}

# --- Notes Master: Date Placeholder 2 (shape index 2) ---
$notesMaster = $p.NotesMaster
$notesMaster.Shapes.Item(2).TextFrame.TextRange.Text = "10/26/22"

# --- Slide 3: recolor two "Straight Arrow Connector" line shapes ---
$slide3 = $p.Slides.Item(3)

# "Straight Arrow Connector 72" (shape 1): line had no fill -> add solid red fill (C00000)
$conn72 = $slide3.Shapes.Item(1)
$conn72.Line.ForeColor.RGB = 192        # RGB(192,0,0) = 0xC00000

# "Straight Arrow Connector 51" (shape 12): solid red fill -> theme accent1 color
$conn51 = $slide3.Shapes.Item(12)
$conn51.Line.ForeColor.ObjectThemeColor = 5   # msoThemeColorAccent1
